$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.103.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.88%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.889.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.33%  '

$ws.Range('E5').Value = '  -1.87%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.24%  '

$ws.Range('E7').Value = '  +0.28%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3171'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.53%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07179'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.86%  '

$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08347'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.54%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7580'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.419'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.99%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.896.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.168.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.78%  '

$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.159'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '251.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.00%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007875'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.79%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.178.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.41%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.936'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1570'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.86%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.301'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.77%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.053'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.480'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.44%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.573'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.538'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.200'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.44%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05354'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7710'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.57%  '

$ws.Range('E37').Value = '  -0.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.733'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.80%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01960'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.62%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.763'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4565'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.099.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.075'
$ws.Range('D43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8757'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.82%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.49'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.89%  '

$ws.Range('E47').Value = '  +0.34%  '

$ws.Range('E48').Value = '  +0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.592'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.76%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.097.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.92%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.576'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.98%  '
